$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- Cells that change from a numeric style to the text placeholder style (General/"0" or "***.*") ---
# Strategy: copy format from an existing placeholder cell of the same kind, then copy its value (shared text).
$src = $ws.Range("D23")
$src.Copy()
$dst = $ws.Range("C15")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

$src = $ws.Range("D23")
$src.Copy()
$dst = $ws.Range("D15")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

$src = $ws.Range("E23")
$src.Copy()
$dst = $ws.Range("E15")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

$src = $ws.Range("D23")
$src.Copy()
$dst = $ws.Range("C23")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

$src = $ws.Range("D23")
$src.Copy()
$dst = $ws.Range("D27")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

$src = $ws.Range("E23")
$src.Copy()
$dst = $ws.Range("E27")
$dst.PasteSpecial(-4122)
$src.Copy()
$dst.PasteSpecial(-4163)

# --- Cells that change from the text placeholder style to a numeric style ---
$src = $ws.Range("C22")
$src.Copy()
$dst = $ws.Range("D22")
$dst.PasteSpecial(-4122)
$dst.Value = 1

$src = $ws.Range("H15")
$src.Copy()
$dst = $ws.Range("E22")
$dst.PasteSpecial(-4122)
$dst.Value = 100

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = -29.411764705882
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -35.593220338983
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -24
$ws.Range("N16").Value = -78.285714285714
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -70
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = -13.157894736842
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 102
$ws.Range("K17").Value = -5.882352941176
$ws.Range("L17").Value = 21.518987341772
$ws.Range("M17").Value = 6.666666666666
$ws.Range("N17").Value = 17.073170731707
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 29.411764705882
$ws.Range("L18").Value = 2.325581395348
$ws.Range("M18").Value = -18.518518518518
$ws.Range("N18").Value = -74.117647058823
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -42.5
$ws.Range("I19").Value = 69
$ws.Range("J19").Value = 94
$ws.Range("K19").Value = -26.595744680851
$ws.Range("L19").Value = -28.125
$ws.Range("M19").Value = 16.949152542372
$ws.Range("N19").Value = -10.389610389610
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 450
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 90
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = 23.529411764705
$ws.Range("L20").Value = -46.153846153846
$ws.Range("M20").Value = 162.5
$ws.Range("N20").Value = -46.153846153846
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 6.451612903225
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -11.965811965812
$ws.Range("I21").Value = 297
$ws.Range("J21").Value = 332
$ws.Range("K21").Value = -10.542168674698
$ws.Range("L21").Value = -17.955801104972
$ws.Range("M21").Value = 9.594095940959
$ws.Range("N21").Value = -50
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 75
$ws.Range("M22").Value = 75
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 30
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = 35.714285714285
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 162
$ws.Range("K24").Value = 23.456790123456
$ws.Range("L24").Value = 23.456790123456
$ws.Range("M24").Value = 100
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 43.75
$ws.Range("I25").Value = 61
$ws.Range("J25").Value = 46
$ws.Range("K25").Value = 32.608695652173
$ws.Range("L25").Value = 19.607843137254
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 25
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -21.818181818181
$ws.Range("I26").Value = 133
$ws.Range("J26").Value = 151
$ws.Range("K26").Value = -11.920529801324
$ws.Range("L26").Value = 29.126213592233
$ws.Range("M26").Value = 15.652173913043
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = -17.647058823529
$ws.Range("L27").Value = 55.555555555555
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 71.428571428571
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 9.523809523809
$ws.Range("L28").Value = 0
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 6
$ws.Range("J30").Value = 5

Write-Output "edit complete"